# Elimna EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" column (E16:E44) listed the billing periods in
# ascending chronological order (1711, 1712, 1801, ... 2003). The data
# was refreshed/resorted so the periods now appear in descending
# (most-recent-first) order: 2003, 2002, 2001, 1912, ... 1712, 1711.
#
# The underlying cells (row/column positions) are unchanged - only the
# text values shown in column E for rows 16 through 44 change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @(
    "2003","2002","2001",
    "1912","1911","1910","1909","1908","1907","1906","1905","1904","1903","1902","1901",
    "1812","1811","1810","1809","1808","1807","1806","1805","1804","1803","1802","1801",
    "1712","1711"
)

$startRow = 16
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
}
